# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for Mango at Terminal Hortofrutícola
# Agro Chillán just before the current row 174, pushing existing rows
# 174-191 down to 175-192 (dimension grows from A1:T191 to A1:T192).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 174, shifting row 174 and everything below
# it down by one (entire-row insert implicitly shifts down).
$ws.Rows.Item(174).Insert()

# Populate the newly inserted row 174 with the new weekly record.
$ws.Cells.Item(174, 1).Value  = 7
$ws.Cells.Item(174, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(174, 3).Value  = "Ñuble"
$ws.Cells.Item(174, 4).Value  = 45194
$ws.Cells.Item(174, 5).Value  = 16
$ws.Cells.Item(174, 6).Value  = "Fruta"
$ws.Cells.Item(174, 7).Value  = 100108
$ws.Cells.Item(174, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(174, 9).Value  = 100108002
$ws.Cells.Item(174, 10).Value = "Mango"
$ws.Cells.Item(174, 11).Value = "Sin especificar"
$ws.Cells.Item(174, 12).Value = "Primera"
$ws.Cells.Item(174, 13).Value = 80
$ws.Cells.Item(174, 14).Value = 10000
$ws.Cells.Item(174, 15).Value = 10000
$ws.Cells.Item(174, 16).Value = 10000
$ws.Cells.Item(174, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(174, 18).Value = "Brasil"
$ws.Cells.Item(174, 19).Value = 2500
$ws.Cells.Item(174, 20).Value = 4
